$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "-"
$ws.Range("D3").Value = "MCT-2A-CAD"
$ws.Range("E3").Value = "MEC-2A-CAD"
$ws.Range("D4").Value = "MCT-2A-CAD"
$ws.Range("E4").Value = "MEC-2A-CAD"
$ws.Range("E6").Value = "MEC-1A-Desenho Técnico"
$ws.Range("E7").Value = "-"
